# Fruta / hortaliza, semanal
# Rotates the weekly price records in rows 4-12 (columns D, L, M, N, O, P, R, S)
# so that the data now reflects the updated weekly reporting order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row for columns D, L, M, N, O, P, R, S
$data = @{
    4  = @{ D = 44294; L = "Primera";  M = 20; N = 225000; O = 230000; P = 227500; R = "Región Metropolitana";   S = 506 }
    5  = @{ D = 44294; L = "Segunda";  M = 16; N = 195000; O = 200000; P = 197500; R = "Región Metropolitana";   S = 439 }
    6  = @{ D = 44309; L = "Especial"; M = 20; N = 305000; O = 310000; P = 307500; R = "Provincia de Cachapoal"; S = 683 }
    7  = @{ D = 44309; L = "Primera";  M = 20; N = 285000; O = 290000; P = 287500; R = "Provincia de Cachapoal"; S = 639 }
    8  = @{ D = 44309; L = "Segunda";  M = 20; N = 255000; O = 260000; P = 257500; R = "Provincia de Cachapoal"; S = 572 }
    9  = @{ D = 44316; L = "Especial"; M = 20; N = 255000; O = 260000; P = 257500; R = "Región de O'Higgins";    S = 572 }
    10 = @{ D = 44316; L = "Primera";  M = 20; N = 225000; O = 230000; P = 227500; R = "Región de O'Higgins";    S = 506 }
    11 = @{ D = 44273; L = "Especial"; M = 10; N = 255000; O = 260000; P = 257500; R = "Región de O'Higgins";    S = 572 }
    12 = @{ D = 44273; L = "Primera";  M = 20; N = 225000; O = 230000; P = 227500; R = "Región de O'Higgins";    S = 506 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("R$row").Value = $vals.R
    $ws.Range("S$row").Value = $vals.S
}
